# Data Contract Updates - commit by 'Lakshmi'
# Updates the "Endpoint Definition" sheet: AccountSummary + Transaction Details
# endpoints get real path params / contracts, and the stray "Admin Approval"
# row is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Endpoint Definition")
$ws.Activate()

# --- Row 5: AccountSummary -----------------------------------------------
$ws.Range("D5").Value = "/api/accountSummary/{customerId}"
$ws.Range("E5").Value = "{`n}"
$ws.Range("F5").Value = "{`n accountNumber:String,`n balance:Double,`n accountType:String`n        }`n"

# --- Row 6: Transaction Details -------------------------------------------
$ws.Range("D6").Value = "/api/transactions/{accountNumber}"
$ws.Range("E6").Value = "{}"
$ws.Range("F6").Value = "{amount:Double,traansactionDate;Date,transactionTime:Time}"

# --- Row 7: remove the old "Admin Approval" endpoint row -------------------
$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()
$ws.Range("D7").ClearContents()

# --- Column D is now wider to fit the longer endpoint paths ---------------
$ws.Columns.Item(4).ColumnWidth = 25.25

# --- Update the view: scroll up a bit and move the selection --------------
$excel.Goto($ws.Range("A3"), $true)
$ws.Range("E4").Select()

Write-Host "Endpoint Definition sheet updated"
